$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Fix the start time on row 42 (was mistakenly entered as 12:52:00 instead of 00:52:00)
$ws.Range("D42").Value = 0.036111111111111115

# Fill in the new activity-log entry on row 43
$ws.Range("B43").Value = 6977
$ws.Range("C43").Value = 43926
$ws.Range("D43").Value = 0.059722222222222225
$ws.Range("E43").Value = 0.06874999999999999
$ws.Range("G43").Value = "Updated Timing waveforms for LogicUnit.vhd"

# Update the active selection to match where the user left off
$ws.Range("C42").Select()
